$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 19328.5
$ws.Range("J3").Value = 19328.5
$ws.Range("L3").Value = 19328.5
$ws.Range("N3").Value = -19556.5
$ws.Range("H19").Value = 165.66667
$ws.Range("I19").Value = 165.66667
$ws.Range("K19").Value = 165.66667
$ws.Range("M19").Value = 9.333329999999989
$ws.Range("H64").Value = 8500
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 8500
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 8500
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -8996
$ws.Range("H67").Value = 8500
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 8500
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 8500
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -10216
$ws.Range("H100").Value = 2799
$ws.Range("I100").Value = 3999.3333
$ws.Range("J100").Value = 998.5
$ws.Range("K100").Value = 3999.3333
$ws.Range("L100").Value = 998.5
$ws.Range("M100").Value = -3458.3333
$ws.Range("N100").Value = -2080.5
$ws.Range("H102").Value = 19328.5
$ws.Range("J102").Value = 19328.5
$ws.Range("L102").Value = 19328.5
$ws.Range("N102").Value = -25818.5
$ws.Range("H132").Value = 31253582
$ws.Range("I132").Value = 35717756
$ws.Range("K132").Value = 107153268
$ws.Range("M132").Value = -107150738
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3499.5
$ws.Range("I45").Value = 3499
$ws.Range("K45").Value = 3499
$ws.Range("M45").Value = -3122
$ws.Range("H61").Value = 2498.5
$ws.Range("J61").Value = 2999
$ws.Range("L61").Value = 2999
$ws.Range("N61").Value = -3423
$ws.Range("H74").Value = 1332.8
$ws.Range("I74").Value = 916
$ws.Range("K74").Value = 916
$ws.Range("M74").Value = -42
$ws.Range("H77").Value = 1332.8
$ws.Range("I77").Value = 916
$ws.Range("K77").Value = 4580
$ws.Range("M77").Value = -212
$ws.Range("H102").Value = 300
$ws.Range("I102").Value = 300
$ws.Range("K102").Value = 300
$ws.Range("M102").Value = 1322
$ws.Range("H132").Value = 4241.55
$ws.Range("I132").Value = 3885.842
$ws.Range("K132").Value = 11657.526
$ws.Range("M132").Value = -9127.526
$ws.Range("H136").Value = 2498.5
$ws.Range("J136").Value = 2999
$ws.Range("L136").Value = 8997
$ws.Range("N136").Value = -14097
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3266.5
$ws.Range("I94").Value = 3266.5
$ws.Range("K94").Value = 3266.5
$ws.Range("M94").Value = -2815.5
$ws.Range("H99").Value = 1490
$ws.Range("I99").Value = 725
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 725
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 773
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 997.5
$ws.Range("I105").Value = 997.5
$ws.Range("K105").Value = 997.5
$ws.Range("M105").Value = 749.5
$ws.Range("H134").Value = 8288.833000000001
$ws.Range("I134").Value = 8666.700000000001
$ws.Range("K134").Value = 26000.1
$ws.Range("M134").Value = -23465.1
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3273.6
$ws.Range("J31").Value = 7770
$ws.Range("L31").Value = 7770
$ws.Range("N31").Value = -8360
$ws.Range("H34").Value = 3273.6
$ws.Range("J34").Value = 7770
$ws.Range("L34").Value = 7770
$ws.Range("N34").Value = -8174
$ws.Range("H58").Value = 2377.5
$ws.Range("I58").Value = 2036.6666
$ws.Range("J58").Value = 3400
$ws.Range("K58").Value = 2036.6666
$ws.Range("L58").Value = 3400
$ws.Range("M58").Value = -1833.6666
$ws.Range("N58").Value = -3806
$ws.Range("H105").Value = 1720.5
$ws.Range("I105").Value = 1720.5
$ws.Range("K105").Value = 1720.5
$ws.Range("M105").Value = 26.5
$ws.Range("H136").Value = 2377.5
$ws.Range("I136").Value = 2036.6666
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 6109.9998
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -3559.9998
$ws.Range("N136").Value = -15300
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 183
$ws.Range("I36").Value = 183
$ws.Range("K36").Value = 549
$ws.Range("M36").Value = -380
$ws.Range("H80").Value = 2499.5
$ws.Range("J80").Value = 2499.5
$ws.Range("L80").Value = 7498.5
$ws.Range("N80").Value = -9370.5
$ws.Range("H83").Value = 2499.5
$ws.Range("J83").Value = 2499.5
$ws.Range("L83").Value = 22495.5
$ws.Range("N83").Value = -31855.5
$ws.Range("H93").Value = 3341
$ws.Range("J93").Value = 3341
$ws.Range("L93").Value = 10023
$ws.Range("N93").Value = -13767
$ws.Range("H109").Value = 2451.5
$ws.Range("I109").Value = 941.8
$ws.Range("K109").Value = 2825.4
$ws.Range("M109").Value = -1785.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 19000000
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H43").Value = 24879.334
$ws.Range("I43").Value = 1749.5
$ws.Range("K43").Value = 1749.5
$ws.Range("M43").Value = -1598.5
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
$ws.Range("H134").Value = 199998
$ws.Range("J134").Value = 199998
$ws.Range("L134").Value = 599994
$ws.Range("N134").Value = -605064
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2944.4167
$ws.Range("I16").Value = 3253.3
$ws.Range("J16").Value = 1400
$ws.Range("K16").Value = 3253.3
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = -3083.3
$ws.Range("N16").Value = -1740
$ws.Range("H69").Value = 47935
$ws.Range("J69").Value = 47935
$ws.Range("L69").Value = 47935
$ws.Range("N69").Value = -49557
$ws.Range("H72").Value = 47935
$ws.Range("J72").Value = 47935
$ws.Range("L72").Value = 143805
$ws.Range("N72").Value = -151917
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H122").Value = 4034.5557
$ws.Range("I122").Value = 3559.4
$ws.Range("K122").Value = 10678.2
$ws.Range("M122").Value = -8228.200000000001
$ws.Range("H132").Value = 13390.25
$ws.Range("I132").Value = 14303
$ws.Range("K132").Value = 42909
$ws.Range("M132").Value = -40379
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2598.3333
$ws.Range("I126").Value = 2598.3333
$ws.Range("K126").Value = 7794.999899999999
$ws.Range("M126").Value = -5324.999899999999
$ws.Range("H132").Value = 1044.75
$ws.Range("I132").Value = 892
$ws.Range("J132").Value = 1197.5
$ws.Range("K132").Value = 2676
$ws.Range("L132").Value = 3592.5
$ws.Range("M132").Value = -146
$ws.Range("N132").Value = -8652.5
